{"js": "// Add a new paragraph style \"CompactList\" (\"Compact List\") that mirrors the\n// existing \"Compact\" style: based on Body Text, quick style, with\n// before/after paragraph spacing of 36 twips (1.8pt).\ncontext.document.addStyle(\"Compact List\", Word.StyleType.paragraph);\nawait context.sync();\n\n// `addStyle`'s returned proxy doesn't reliably resolve right away in all\n// hosts, so re-acquire the style by name before configuring it.\nconst styles = context.document.getStyles();\nconst style = styles.getByNameOrNullObject(\"Compact List\");\nawait context.sync();\n\n// Base it on the \"Body Text\" style (styleId \"BodyText\") \u2014 same as \"Compact\".\nstyle.baseStyle = \"BodyText\";\n// Mark it as a Quick Style (<w:qFormat/>), same as \"Compact\".\nstyle.quickStyle = true;\n// Paragraph spacing before/after = 36 dxa = 1.8 pt, same as \"Compact\".\nstyle.paragraphFormat.spaceBefore = 1.8;\nstyle.paragraphFormat.spaceAfter = 1.8;\n\nawait context.sync();\n", "ps1": "# Add a new paragraph style \"CompactList\" (\"Compact List\") that mirrors the\n# existing \"Compact\" style: based on Body Text, quick style, with\n# before/after paragraph spacing of 36 twips (1.8pt).\n$d = $word.ActiveDocument\n\n# wdStyleTypeParagraph = 1\n$style = $d.Styles.Add(\"Compact List\", 1)\n\n# Base it on the \"Body Text\" style (styleId \"BodyText\") \u2014 same as \"Compact\".\n# Use the style id string so w:basedOn stores the id, not the display name.\n$style.BaseStyle = \"BodyText\"\n\n# Mark it as a Quick Style (<w:qFormat/>), same as \"Compact\".\n$style.QuickStyle = $true\n\n# Paragraph spacing before/after = 36 dxa = 1.8 pt, same as \"Compact\".\n$style.ParagraphFormat.SpaceBefore = 1.8\n$style.ParagraphFormat.SpaceAfter = 1.8\n"}
